$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert two new rows inside the leave-card table to record the new leave
#    events (SL(1-0-0), SP(3-0-0), FL(2-0-0)) taken around 12/2023.
#    Row 80 (new) -> pushes the existing monthly-accrual row down to 81.
#    Row 82 (new, inserted after the first shift) -> pushes everything else
#    down so the table grows from A8:K143 to A8:K145.
# ---------------------------------------------------------------------------
$ws.Rows.Item(80).Insert()
$ws.Rows.Item(82).Insert()

# Re-attach the table to the now-larger range (Excel normally grows the
# table automatically when a row is inserted inside it; make sure it covers
# the two appended rows at the bottom as well).
$tbl.Resize($ws.Range("A8:K145"))

# ---------------------------------------------------------------------------
# 2. The freshly inserted rows (80 and 82) come in with generic/default
#    formatting. Restore the table's normal data-row look (borders, number
#    formats, etc.) by copying the formatting from an existing data row, then
#    restore the table's calculated "EARNED " column formula.
# ---------------------------------------------------------------------------
$earnedFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$ws.Range("A79:K79").Copy()
$ws.Range("A80:K80").PasteSpecial(-4122)
$ws.Range("G80").Formula = $earnedFormula

$ws.Range("A79:K79").Copy()
$ws.Range("A82:K82").PasteSpecial(-4122)
$ws.Range("G82").Formula = $earnedFormula

$ws.Application.CutCopyMode = $false

# The two rows that used to sit at the very bottom of the table (old 143,
# the specially-styled closing row) shifted down to 144/145 along with the
# insert above; their "EARNED " formula needs to keep using the table
# structured reference (the plain insert collapses it to [@EARNED]).
$ws.Range("G144").Formula = $earnedFormula
$ws.Range("G145").Formula = $earnedFormula

# ---------------------------------------------------------------------------
# 3. Fill in the actual leave-card data for the new / shifted rows.
# ---------------------------------------------------------------------------

# Monthly SL credit entries that became due after the new rows were added.
$ws.Range("C75").Value = 1.25
$ws.Range("C76").Value = 1.25
$ws.Range("C78").Value = 1.25

# Row 77: SL(1-0-0) used, approved 9/26/2023.
$ws.Range("B77").Value = "SL(1-0-0)"
$ws.Range("C77").Value = 1.25
$ws.Range("H77").Value = 1
$ws.Range("K72").Copy()
$ws.Range("K77").PasteSpecial(-4122)
$ws.Range("K77").Value = 45195

# Row 79: SL(1-0-0) used, approved 11/3/2023.
$ws.Range("B79").Value = "SL(1-0-0)"
$ws.Range("C79").Value = 1.25
$ws.Range("H79").Value = 1
$ws.Range("K72").Copy()
$ws.Range("K79").PasteSpecial(-4122)
$ws.Range("K79").Value = 45233

# Row 80 (new): SL(1-0-0) used, approved 11/8/2023.
$ws.Range("B80").Value = "SL(1-0-0)"
$ws.Range("H80").Value = 1
$ws.Range("K72").Copy()
$ws.Range("K80").PasteSpecial(-4122)
$ws.Range("K80").Value = 45238

# Row 81 (was the monthly-accrual row for 12/1/2023): SP(3-0-0) used,
# covering 12/20-22/2023.
$ws.Range("B81").Value = "SP(3-0-0)"
$ws.Range("K81").Value = "12/20-22/2023"

# Row 82 (new): FL(2-0-0) used, covering 12/18,19/2023.
$ws.Range("B82").Value = "FL(2-0-0)"
$ws.Range("D82").Value = 2
$ws.Range("K82").Value = "12/18,19/2023"

$ws.Application.CutCopyMode = $false
